$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 4: PWM Fans x2 price updated
# ---------------------------------------------------------------
$ws.Range("F4").Value = 35.770000000000003

# ---------------------------------------------------------------
# Row 10 (PCB row) content changes
# ---------------------------------------------------------------
$ws.Range("A10").Value = "PCB"
$ws.Range("B10").Value = "Elecrow"
$ws.Range("C10").Value = "Elecrow"
$ws.Range("F10").Value = 4.99
$ws.Range("G10").Value = 1
$ws.Range("H10").Formula = "=F10*G10"
$ws.Range("I10").Value = "26.99  Fed Ex"
$ws.Range("J10").Value = "To order from Elecrow: Default settings should be correct. But just to be sure, the settings are:`n1) Layers: 2`n2) Dimensions: 100 x 94mm (same price as 100x100)`n3) Quantity 5`n4) Different pcb design 1`n5) PCB Thickness: 1.6mm`n6) color: green`n7) surface finish: Hasl`n8) castellated hole: no`n9) Coper weight: 1oz`n10) PCB stencil: No`nShipping: DHL"

$ws.Rows.Item(10).RowHeight = 225

# New hyperlink on B10 pointing at the Elecrow PCB service
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.elecrow.com/pcb-prototyping.html") | Out-Null

# ---------------------------------------------------------------
# Row 11 (Totals row)
# ---------------------------------------------------------------
$ws.Rows.Item(11).RowHeight = 23.25
$ws.Range("I11").Value = 35.979999999999997
$ws.Range("J11").Formula = "=H11+I11"
$ws.Range("J11").Font.Bold = $true
$ws.Range("J11").Font.Size = 18

Write-Output "edit complete"
